$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections in rows 2-25 (data unaffected by the row reshuffle below) ---

# RM 8 (row 3): "D" column (sheet col E) had no value -> now -5.7
$ws.Cells.Item(3, 5).Value = -5.7

# RM 9 (row 4): F value 17.97 removed (no longer measured)
$ws.Cells.Item(4, 6).ClearContents()

# RM 14 (row 5): E value -5 removed (no longer measured)
$ws.Cells.Item(5, 5).ClearContents()

# RM 42 (row 9): F value now 17.26
$ws.Cells.Item(9, 6).Value = 17.26

# RM 52 a (row 10): F value now 16.43
$ws.Cells.Item(10, 6).Value = 16.43

# RM 116 (row 17): F value 17.78 removed
$ws.Cells.Item(17, 6).ClearContents()

# RM 120 (row 18): F value 18.35 removed
$ws.Cells.Item(18, 6).ClearContents()

# RM 135 (row 21): E value now -8.699999999999999
$ws.Cells.Item(21, 5).Value = -8.699999999999999

# RM 140 (row 23): E value -7 removed
$ws.Cells.Item(23, 5).ClearContents()

# --- Remove the "RM 232" and "SC 92" rows, shifting the remaining SC rows up ---
# Delete the lower row first so the row index of the upper one doesn't move.
$ws.Rows(28).Delete()   # "SC 92" row (originally row 28)
$ws.Rows(26).Delete()   # "RM 232" row (originally row 26)

# --- Post-shift fix: "SC 193" (now row 32) gained an E value ---
$ws.Cells.Item(32, 5).Value = -6.4
